# Fix typo'd / stray-character state names left over from test data, and
# remove the trailing non-breaking space from "Tennessee" on the totals row.
# (commit: "Removed test scripts and functions")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A51").Value = "Wyoming"    # was "Wyominga"
$ws.Range("A13").Value = "Idaho"      # was "Idahoa"
$ws.Range("A17").Value = "Kansas"     # was "Kansasa"
$ws.Range("A26").Value = "Missouri"   # was "Missouria"
$ws.Range("A52").Value = "Tennessee"  # was "Tennessee " (trailing NBSP)

# Update the active selection to match the saved view (scrolled to the
# bottom of the state list).
$ws.Activate()
$ws.Range("A52").Select()
